$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1657.7441
$ws.Range("I40").Value = 1503.8462
$ws.Range("J40").Value = 1893.1177
$ws.Range("K40").Value = 1503.8462
$ws.Range("L40").Value = 1893.1177
$ws.Range("M40").Value = -1328.8462
$ws.Range("N40").Value = -2243.1177
$ws.Range("H53").Value = 128.63333
$ws.Range("I53").Value = 89.52941
$ws.Range("K53").Value = 89.52941
$ws.Range("M53").Value = 547.47059
$ws.Range("H55").Value = 258.7143
$ws.Range("I55").Value = 192.2
$ws.Range("J55").Value = 425
$ws.Range("K55").Value = 192.2
$ws.Range("L55").Value = 425
$ws.Range("M55").Value = 21.80000000000001
$ws.Range("N55").Value = -853
$ws.Range("H107").Value = 979.3158
$ws.Range("I107").Value = 926.8182
$ws.Range("J107").Value = 1051.5
$ws.Range("K107").Value = 926.8182
$ws.Range("L107").Value = 1051.5
$ws.Range("M107").Value = 993.1818
$ws.Range("N107").Value = -4891.5
$ws.Range("H116").Value = 99662.37
$ws.Range("I116").Value = 135325.75
$ws.Range("J116").Value = 4560
$ws.Range("K116").Value = 135325.75
$ws.Range("L116").Value = 4560
$ws.Range("M116").Value = -131883.75
$ws.Range("N116").Value = -11444
$ws.Range("H138").Value = 1982.7764
$ws.Range("I138").Value = 868.25
$ws.Range("J138").Value = 3221.139
$ws.Range("K138").Value = 2604.75
$ws.Range("L138").Value = 9663.417000000001
$ws.Range("M138").Value = 2535.25
$ws.Range("N138").Value = -19943.417

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1160178.8
$ws.Range("I32").Value = 1324984.8
$ws.Range("J32").Value = 6536.4443
$ws.Range("K32").Value = 1324984.8
$ws.Range("L32").Value = 6536.4443
$ws.Range("M32").Value = -1324697.8
$ws.Range("N32").Value = -7110.4443
$ws.Range("H132").Value = 15938.918
$ws.Range("I132").Value = 21869.8
$ws.Range("J132").Value = 3045.6956
$ws.Range("K132").Value = 65609.39999999999
$ws.Range("L132").Value = 9137.086800000001
$ws.Range("M132").Value = -63079.39999999999
$ws.Range("N132").Value = -14197.0868

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2189.5715
$ws.Range("I20").Value = 1999.5
$ws.Range("J20").Value = 2265.6
$ws.Range("K20").Value = 1999.5
$ws.Range("L20").Value = 2265.6
$ws.Range("M20").Value = -1752.5
$ws.Range("N20").Value = -2759.6
$ws.Range("H134").Value = 3109.2273
$ws.Range("I134").Value = 2860.7878
$ws.Range("K134").Value = 8582.3634
$ws.Range("M134").Value = -6047.3634

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2414.4126
$ws.Range("I31").Value = 1664.5581
$ws.Range("J31").Value = 4026.6
$ws.Range("K31").Value = 1664.5581
$ws.Range("L31").Value = 4026.6
$ws.Range("M31").Value = -1369.5581
$ws.Range("N31").Value = -4616.6
$ws.Range("H34").Value = 2414.4126
$ws.Range("I34").Value = 1664.5581
$ws.Range("J34").Value = 4026.6
$ws.Range("K34").Value = 1664.5581
$ws.Range("L34").Value = 4026.6
$ws.Range("M34").Value = -1462.5581
$ws.Range("N34").Value = -4430.6
$ws.Range("H58").Value = 1255.4348
$ws.Range("I58").Value = 704.1515000000001
$ws.Range("J58").Value = 2654.8462
$ws.Range("K58").Value = 704.1515000000001
$ws.Range("L58").Value = 2654.8462
$ws.Range("M58").Value = -501.1515000000001
$ws.Range("N58").Value = -3060.8462
$ws.Range("H64").Value = 39000
$ws.Range("J64").Value = 39000
$ws.Range("L64").Value = 39000
$ws.Range("N64").Value = -39496
$ws.Range("H67").Value = 39000
$ws.Range("J67").Value = 39000
$ws.Range("L67").Value = 39000
$ws.Range("N67").Value = -40716
$ws.Range("H122").Value = 1156.0625
$ws.Range("I122").Value = 998
$ws.Range("J122").Value = 1314.125
$ws.Range("K122").Value = 2994
$ws.Range("L122").Value = 3942.375
$ws.Range("M122").Value = -544
$ws.Range("N122").Value = -8842.375
$ws.Range("H134").Value = 2030.1111
$ws.Range("I134").Value = 1260.3846
$ws.Range("J134").Value = 2744.8572
$ws.Range("K134").Value = 3781.1538
$ws.Range("L134").Value = 8234.571599999999
$ws.Range("M134").Value = -1246.1538
$ws.Range("N134").Value = -13304.5716
$ws.Range("H136").Value = 1255.4348
$ws.Range("I136").Value = 704.1515000000001
$ws.Range("J136").Value = 2654.8462
$ws.Range("K136").Value = 2112.4545
$ws.Range("L136").Value = 7964.5386
$ws.Range("M136").Value = 437.5454999999997
$ws.Range("N136").Value = -13064.5386

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 3588.8333
$ws.Range("I56").Value = 3588.8333
$ws.Range("K56").Value = 3588.8333
$ws.Range("M56").Value = -3058.8333
$ws.Range("H132").Value = 7438.3076
$ws.Range("I132").Value = 4279.6
$ws.Range("J132").Value = 9412.5
$ws.Range("K132").Value = 38516.4
$ws.Range("L132").Value = 84712.5
$ws.Range("M132").Value = -35986.4
$ws.Range("N132").Value = -89772.5
$ws.Range("H137").Value = 2879.24
$ws.Range("I137").Value = 1353.3334
$ws.Range("J137").Value = 5168.1
$ws.Range("K137").Value = 4060.0002
$ws.Range("L137").Value = 15504.3
$ws.Range("M137").Value = 1039.9998
$ws.Range("N137").Value = -25704.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1486.4445
$ws.Range("I122").Value = 1404.421
$ws.Range("J122").Value = 1681.25
$ws.Range("K122").Value = 4213.263
$ws.Range("L122").Value = 5043.75
$ws.Range("M122").Value = -1763.263
$ws.Range("N122").Value = -9943.75
$ws.Range("H126").Value = 2827.7896
$ws.Range("I126").Value = 2256.8635
$ws.Range("J126").Value = 3612.8125
$ws.Range("K126").Value = 6770.5905
$ws.Range("L126").Value = 10838.4375
$ws.Range("M126").Value = -4300.5905
$ws.Range("N126").Value = -15778.4375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1623.6666
$ws.Range("I93").Value = 1448.4
$ws.Range("J93").Value = 2500
$ws.Range("K93").Value = 1448.4
$ws.Range("L93").Value = 2500
$ws.Range("M93").Value = -200.4000000000001
$ws.Range("N93").Value = -4996

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2328.875
$ws.Range("I132").Value = 1674.5238
$ws.Range("J132").Value = 3578.0908
$ws.Range("K132").Value = 5023.5714
$ws.Range("L132").Value = 10734.2724
$ws.Range("M132").Value = -2493.5714
$ws.Range("N132").Value = -15794.2724

Write-Host "Updated leve profit values across 8 sheets."